$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates ---
$ws.Range("B1").Value2 = "desc"
$ws.Range("J1").Value2 = "accuracy (%)"

# --- Existing rows 2-4: rename "yes" -> "bsif", recompute accuracy as percentage ---
$ws.Range("B2").Value2 = "bsif"
$ws.Range("J2").Value2 = 70.76

$ws.Range("B3").Value2 = "bsif"
$ws.Range("J3").Value2 = 70.76

$ws.Range("B4").Value2 = "bsif"
$ws.Range("J4").Value2 = 70.76

# --- New rows 5-8 ---
$ws.Range("A5").Value2 = "simple_test"
$ws.Range("B5").Value2 = "gray"
$ws.Range("C5").Value2 = "-"
$ws.Range("D5").Value2 = "-"
$ws.Range("E5").Value2 = "handout 0.9"
$ws.Range("F5").Value2 = "10x32"
$ws.Range("G5").Value2 = 5.16
$ws.Range("H5").Value2 = "knn mode"
$ws.Range("I5").Value2 = "left"
$ws.Range("J5").Value2 = 60.67

$ws.Range("A6").Value2 = "simple_test"
$ws.Range("B6").Value2 = "gray"
$ws.Range("C6").Value2 = "-"
$ws.Range("D6").Value2 = "-"
$ws.Range("E6").Value2 = "handout 0.9"
$ws.Range("F6").Value2 = "10x32"
$ws.Range("G6").Value2 = 5.58
$ws.Range("H6").Value2 = "knn mode"
$ws.Range("I6").Value2 = "left"
$ws.Range("J6").Value2 = 64

$ws.Range("A7").Value2 = "simple_test"
$ws.Range("B7").Value2 = "bsif"
$ws.Range("C7").Value2 = 15
$ws.Range("D7").Value2 = 11
$ws.Range("E7").Value2 = "handout 0.9"
$ws.Range("F7").Value2 = "10x32"
$ws.Range("G7").Value2 = 8.94
$ws.Range("H7").Value2 = "knn mode"
$ws.Range("I7").Value2 = "left"
$ws.Range("J7").Value2 = 72

$ws.Range("A8").Value2 = "simple_test"
$ws.Range("B8").Value2 = "bsif"
$ws.Range("C8").Value2 = 5
$ws.Range("D8").Value2 = 11
$ws.Range("E8").Value2 = "handout 0.9"
$ws.Range("F8").Value2 = "10x32"
$ws.Range("G8").Value2 = 8.94
$ws.Range("H8").Value2 = "knn mode"
$ws.Range("I8").Value2 = "left"
$ws.Range("J8").Value2 = 72

# --- Match final selection state from the saved workbook ---
$ws.Range("J8").Select() | Out-Null
